$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "72.921.21"
$ws.Range("E2").Value = "  +2.02%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.988.11"
$ws.Range("E3").Value = "  +0.14%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "593.08"
$ws.Range("E5").Value = "  +9.59%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "161.27"
$ws.Range("E6").Value = "  +8.21%  "
$ws.Range("E7").Value = "  -0.60%  "
$ws.Range("E8").Value = "  -0.13%  "
$ws.Range("E9").Value = "  +0.48%  "
$ws.Range("E10").Value = "  +0.62%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.97"
$ws.Range("E11").Value = "  -4.35%  "
$ws.Range("E12").Value = "  -0.25%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.97"
$ws.Range("E13").Value = "  +1.90%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.620.37"
$ws.Range("E14").Value = "  +0.09%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.995.41"
$ws.Range("E15").Value = "  +0.49%  "
$ws.Range("E16").Value = "  +8.33%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.15"
$ws.Range("E17").Value = "  +1.20%  "
$ws.Range("E18").Value = "  -0.92%  "
$ws.Range("E19").Value = "  +0.28%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "72.558.35"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "433.21"
$ws.Range("E21").Value = "  +0.98%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.79"
$ws.Range("E22").Value = "  +13.48%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "96.20"
$ws.Range("E23").Value = "  -1.31%  "
$ws.Range("E24").Value = "  -4.65%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "14.19"
$ws.Range("E25").Value = "  -2.49%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.39"
$ws.Range("E26").Value = "  +17.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.31"
$ws.Range("E27").Value = "  -1.85%  "
$ws.Range("E28").Value = "  +0.56%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.49"
$ws.Range("E29").Value = "  -2.60%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.37"
$ws.Range("E30").Value = "  -0.85%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.90"
$ws.Range("E31").Value = "  +2.85%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "13.72"
$ws.Range("E32").Value = "  +2.18%  "
$ws.Range("E33").Value = "  -0.24%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "48.86"
$ws.Range("E34").Value = "  -4.66%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "671.32"
$ws.Range("E35").Value = "  -2.28%  "
$ws.Range("E36").Value = "  +7.36%  "
$ws.Range("B37").Value = "PEPE"
$ws.Range("C37").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D37").Value = "0.0₃0882"
$ws.Range("E37").Value = "  +7.30%  "
$ws.Range("B38").Value = "TheGraph"
$ws.Range("C38").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.436"
$ws.Range("E38").Value = "  -0.53%  "
$ws.Range("E39").Value = "  -3.17%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.997"
$ws.Range("E40").Value = "  +0.04%  "
$ws.Range("E41").Value = "  -1.91%  "
$ws.Range("E42").Value = "  +1.83%  "
$ws.Range("E43").Value = "  +0.03%  "
$ws.Range("E44").Value = "  +0.85%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.69"
$ws.Range("E45").Value = "  +10.16%  "
$ws.Range("E46").Value = "  +0.26%  "
$ws.Range("E47").Value = "  +2.25%  "
$ws.Range("E48").Value = "  -4.05%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.841.95"
$ws.Range("E49").Value = "  +3.66%  "
$ws.Range("E50").Value = "  +0.30%  "
$ws.Range("E51").Value = "  +3.95%  "
